# Refresh latest market-price snapshot values (currentAveragePrice / NQ / HQ
# columns and their dependent Leve profit calculations) across all profession
# sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4199.8184
$ws.Range("I40").Value = 3400.3333
$ws.Range("J40").Value = 4499.625
$ws.Range("K40").Value = 3400.3333
$ws.Range("L40").Value = 4499.625
$ws.Range("M40").Value = -3225.3333
$ws.Range("N40").Value = -4849.625
$ws.Range("H64").Value = 40006796
$ws.Range("I64").Value = 55562436
$ws.Range("J64").Value = 6571.4287
$ws.Range("K64").Value = 55562436
$ws.Range("L64").Value = 6571.4287
$ws.Range("M64").Value = -55562188
$ws.Range("N64").Value = -7067.4287
$ws.Range("H67").Value = 40006796
$ws.Range("I67").Value = 55562436
$ws.Range("J67").Value = 6571.4287
$ws.Range("K67").Value = 55562436
$ws.Range("L67").Value = 6571.4287
$ws.Range("M67").Value = -55561578
$ws.Range("N67").Value = -8287.4287
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H129").Value = 1274.35
$ws.Range("I129").Value = 670.0909
$ws.Range("K129").Value = 2010.2727
$ws.Range("M129").Value = 2989.7273
$ws.Range("H135").Value = 834664.3
$ws.Range("I135").Value = 910457.4399999999
$ws.Range("K135").Value = 8194116.959999999
$ws.Range("M135").Value = -8191581.959999999
$ws.Range("H137").Value = 5667.3706
$ws.Range("I137").Value = 4083.8462
$ws.Range("J137").Value = 7137.7856
$ws.Range("K137").Value = 12251.5386
$ws.Range("L137").Value = 21413.3568
$ws.Range("M137").Value = -9701.5386
$ws.Range("N137").Value = -26513.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1815083
$ws.Range("I32").Value = 1815083
$ws.Range("K32").Value = 1815083
$ws.Range("M32").Value = -1814796
$ws.Range("H61").Value = 66675396
$ws.Range("I61").Value = 1275.5714
$ws.Range("J61").Value = 125015250
$ws.Range("K61").Value = 1275.5714
$ws.Range("L61").Value = 125015250
$ws.Range("M61").Value = -1063.5714
$ws.Range("N61").Value = -125015674
$ws.Range("H74").Value = 103907
$ws.Range("J74").Value = 5137.636
$ws.Range("L74").Value = 5137.636
$ws.Range("N74").Value = -6885.636
$ws.Range("H77").Value = 103907
$ws.Range("J77").Value = 5137.636
$ws.Range("L77").Value = 25688.18
$ws.Range("N77").Value = -34424.18
$ws.Range("H106").Value = 34150.285
$ws.Range("J106").Value = 39410.6
$ws.Range("L106").Value = 39410.6
$ws.Range("N106").Value = -41934.6
$ws.Range("H136").Value = 66675396
$ws.Range("I136").Value = 1275.5714
$ws.Range("J136").Value = 125015250
$ws.Range("K136").Value = 3826.7142
$ws.Range("L136").Value = 375045750
$ws.Range("M136").Value = -1276.7142
$ws.Range("N136").Value = -375050850

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 937.3889
$ws.Range("I22").Value = 853.9231
$ws.Range("J22").Value = 1154.4
$ws.Range("K22").Value = 853.9231
$ws.Range("L22").Value = 1154.4
$ws.Range("M22").Value = -503.9231
$ws.Range("N22").Value = -1854.4
$ws.Range("H31").Value = 5910.5894
$ws.Range("I31").Value = 2342.1738
$ws.Range("K31").Value = 2342.1738
$ws.Range("M31").Value = -2047.1738
$ws.Range("H34").Value = 5910.5894
$ws.Range("I34").Value = 2342.1738
$ws.Range("K34").Value = 2342.1738
$ws.Range("M34").Value = -2140.1738
$ws.Range("H58").Value = 9039.32
$ws.Range("I58").Value = 1772
$ws.Range("K58").Value = 1772
$ws.Range("M58").Value = -1569
$ws.Range("H134").Value = 5029.3335
$ws.Range("I134").Value = 2280.2068
$ws.Range("J134").Value = 10012.125
$ws.Range("K134").Value = 6840.6204
$ws.Range("L134").Value = 30036.375
$ws.Range("M134").Value = -4305.6204
$ws.Range("N134").Value = -35106.375
$ws.Range("H136").Value = 9039.32
$ws.Range("I136").Value = 1772
$ws.Range("K136").Value = 5316
$ws.Range("M136").Value = -2766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 401309.8
$ws.Range("I9").Value = 2275
$ws.Range("J9").Value = 667333
$ws.Range("K9").Value = 6825
$ws.Range("L9").Value = 2001999
$ws.Range("M9").Value = -6601
$ws.Range("N9").Value = -2002447
$ws.Range("H10").Value = 468.22223
$ws.Range("I10").Value = 401.75
$ws.Range("K10").Value = 1205.25
$ws.Range("M10").Value = -1066.25
$ws.Range("H15").Value = 1300.9231
$ws.Range("J15").Value = 1599.9
$ws.Range("L15").Value = 4799.700000000001
$ws.Range("N15").Value = -5079.700000000001
$ws.Range("H121").Value = 2942649.5
$ws.Range("I121").Value = 1388.909
$ws.Range("J121").Value = 8334960.5
$ws.Range("K121").Value = 4166.727000000001
$ws.Range("L121").Value = 25004881.5
$ws.Range("M121").Value = -2856.727000000001
$ws.Range("N121").Value = -25007501.5
$ws.Range("H129").Value = 23882062
$ws.Range("J129").Value = 55723640
$ws.Range("L129").Value = 167170920
$ws.Range("N129").Value = -167180920
$ws.Range("H131").Value = 2396.1707
$ws.Range("I131").Value = 1493.6154
$ws.Range("J131").Value = 2815.2144
$ws.Range("K131").Value = 4480.8462
$ws.Range("L131").Value = 8445.643199999999
$ws.Range("M131").Value = 559.1538
$ws.Range("N131").Value = -18525.6432
$ws.Range("H137").Value = 183586.45
$ws.Range("I137").Value = 112602.11
$ws.Range("J137").Value = 503016
$ws.Range("K137").Value = 337806.33
$ws.Range("L137").Value = 1509048
$ws.Range("M137").Value = -332706.33
$ws.Range("N137").Value = -1519248

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3094.1304
$ws.Range("I102").Value = 3168.3333
$ws.Range("J102").Value = 2827
$ws.Range("K102").Value = 3168.3333
$ws.Range("L102").Value = 2827
$ws.Range("M102").Value = -1546.3333
$ws.Range("N102").Value = -6071

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3906.4443
$ws.Range("I68").Value = 3894.75
$ws.Range("K68").Value = 3894.75
$ws.Range("M68").Value = -3145.75
$ws.Range("H71").Value = 3906.4443
$ws.Range("I71").Value = 3894.75
$ws.Range("K71").Value = 19473.75
$ws.Range("M71").Value = -15729.75
$ws.Range("H82").Value = 529287.75
$ws.Range("I82").Value = 1113653.2
$ws.Range("J82").Value = 3358.8
$ws.Range("K82").Value = 1113653.2
$ws.Range("L82").Value = 3358.8
$ws.Range("M82").Value = -1113292.2
$ws.Range("N82").Value = -4080.8
$ws.Range("H85").Value = 529287.75
$ws.Range("I85").Value = 1113653.2
$ws.Range("J85").Value = 3358.8
$ws.Range("K85").Value = 1113653.2
$ws.Range("L85").Value = 3358.8
$ws.Range("M85").Value = -1112405.2
$ws.Range("N85").Value = -5854.8
$ws.Range("H93").Value = 2595.0557
$ws.Range("I93").Value = 2624
$ws.Range("J93").Value = 2519.8
$ws.Range("K93").Value = 2624
$ws.Range("L93").Value = 2519.8
$ws.Range("M93").Value = -1376
$ws.Range("N93").Value = -5015.8
$ws.Range("H100").Value = 3807
$ws.Range("I100").Value = 1523.4286
$ws.Range("K100").Value = 1523.4286
$ws.Range("M100").Value = -982.4286
$ws.Range("H132").Value = 9443080
$ws.Range("I132").Value = 19233450
$ws.Range("J132").Value = 15315.963
$ws.Range("K132").Value = 57700350
$ws.Range("L132").Value = 45947.889
$ws.Range("M132").Value = -57697820
$ws.Range("N132").Value = -51007.889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2220.1428
$ws.Range("I4").Value = 88
$ws.Range("J4").Value = 7550.5
$ws.Range("K4").Value = 88
$ws.Range("L4").Value = 7550.5
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = -7776.5
$ws.Range("H122").Value = 9886776
$ws.Range("I122").Value = 18005596
$ws.Range("J122").Value = 2995.1738
$ws.Range("K122").Value = 54016788
$ws.Range("L122").Value = 8985.5214
$ws.Range("M122").Value = -54014338
$ws.Range("N122").Value = -13885.5214
